# Generate Report for Handoff
# - Priority for the 3b511c5a... (zh-cn / de-de) rows flips from "low" to "ht"
# - The zh-cn handoff for that same file gets a fresh "Latest Handoff Datetime"
# - The shared "Latest HO Xliff Generate Date" timestamp (Overview + de-de) is refreshed

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-08-28 14:30:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-08-28 14:31:03"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-08-28 14:31:03"
